$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as Excel serial numbers, matching column A date format)
$newRows = @(
    @{ Row = 234; A = 44308; B = 2; C = 11; D = 118.4834123222749 },
    @{ Row = 235; A = 44309; B = 1; C = 9;  D = 96.9409737182249 },
    @{ Row = 236; A = 44310; B = 1; C = 9;  D = 96.9409737182249 },
    @{ Row = 237; A = 44311; B = 3; C = 11; D = 118.4834123222749 },
    @{ Row = 238; A = 44312; B = 2; C = 11; D = 118.4834123222749 }
)

# Use the last existing data row (233) as the style template for the new rows
$templateRow = 233

foreach ($item in $newRows) {
    $r = $item.Row

    # Copy formatting (style) from the template cell in column A only, since that
    # column carries the date number format / border / alignment. Columns B:D keep
    # the workbook's default (unstyled) formatting, matching the existing rows.
    $ws.Range("A$templateRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}

$excel.CutCopyMode = 0
